# Revise immune cell markers
# Append 5 new marker rows for "CD8 T cells" to Sheet1 (rows 36-40)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newMarkers = @("ZFP36", "IFIT1", "GZMK", "IFNG", "LAG3")

$row = 36
foreach ($marker in $newMarkers) {
    $ws.Cells.Item($row, 1).Value = "CD8 T cells"
    $ws.Cells.Item($row, 2).Value = $marker
    $row = $row + 1
}

# Update the selected cell / view to roughly match the final state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D36").Select()
